$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8 (Development, date 2022-01-04)
$ws.Range("A8").Value = 44565
$ws.Range("A8").NumberFormat = "mmm-dd-yy"
$ws.Range("B8").Value = "Development"
$ws.Range("C8").Value = 119
$ws.Range("D8").Value = 117
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "After execution all test cases pass"
$ws.Range("G8").Value = "Test cases iitially fail because of page load affected by network"

# Add new row 9 (Production, date 2022-01-04)
$ws.Range("A9").Value = 44565
$ws.Range("A9").NumberFormat = "mmm-dd-yy"
$ws.Range("B9").Value = "Production"
$ws.Range("C9").Value = 134
$ws.Range("D9").Value = 129
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "After execution all test cases pass"
$ws.Range("G9").Value = "Test cases initially fail because of page load affected by network"

# Apply styling to match existing rows (row height, centered alignment, wrap text)
$ws.Range("A8:G9").RowHeight = 75
$ws.Range("C8:E9").HorizontalAlignment = -4108
$ws.Range("C8:E9").VerticalAlignment = -4108
$ws.Range("F8:G9").WrapText = $true

$ws.Range("A1").Select()
$ws.Range("H9").Select()

$wv = $excel.ActiveWindow
$wv.ScrollRow = 7

$excel.ActiveWindow.DisplayHorizontalScrollBar = $true
